$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.660.43'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.353.82'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.19'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.42'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -7.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -3.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.16'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0921'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.43'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -4.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.997'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -4.30%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.31'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.713.03'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.362.33'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.97'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +9.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.679.62'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000106'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.04'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.70'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +7.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '264.34'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.30'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -10.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.01'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +8.12%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.43'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -5.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.84'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.45'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0900'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -3.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.18'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -10.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.06'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.55'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -8.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0357'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -5.26%  '
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.90'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +4.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -9.29%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.238'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.63'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -4.19%  '
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '119.79'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +7.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.46'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +21.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.88'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -7.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.52'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.22'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -1.81%  '

# Rows 50-51: Cronos / TrustWalletToken swapped positions with updated data
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +0.09%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.26'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -4.31%  '

